$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with merged values across all former rows (2-10)
$ws.Range("A2").Value = "merged_06GE2U92FA_DTPQ67872X_G8WP29EHC6_HYT3MYM7CY_IMFA21HLV3_JU9OS20S99_PPB56V08LB_R5X6KPETN3_RX38XS00QN"
$ws.Range("B2").Value = "merged_L_M_S_XL_XS_XXL_XXS_onesize"
$ws.Range("C2").Value = "merged_Amsterdam_Berlin_Brussels_Copenhagen_Helsinki_Madrid_Paris_Platform_Rome_Stockholm_Webshop"
$ws.Range("D2").Value = 2024
$ws.Range("E2").Value = 12524

# Delete rows 3 through 10 (the rest of the original data that got merged into row 2)
$ws.Range("A3:A10").EntireRow.Delete()
